$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1143.5
$ws.Range("I96").Value = 789.6
$ws.Range("J96").Value = 1733.3334
$ws.Range("K96").Value = 2368.8
$ws.Range("L96").Value = 5200.0002
$ws.Range("M96").Value = -995.8000000000002
$ws.Range("N96").Value = -7946.0002

$ws.Range("H112").Value = 1825
$ws.Range("J112").Value = 2036.6666
$ws.Range("L112").Value = 6109.9998
$ws.Range("N112").Value = -8325.9998

$ws.Range("H121").Value = 628
$ws.Range("J121").Value = 638.5
$ws.Range("L121").Value = 1915.5
$ws.Range("N121").Value = -5409.5

$ws.Range("H128").Value = 35748
$ws.Range("J128").Value = 35748
$ws.Range("L128").Value = 35748
$ws.Range("N128").Value = -45708

$ws.Range("H138").Value = 3168.77
$ws.Range("I138").Value = 1077.7778
$ws.Range("J138").Value = 4344.953
$ws.Range("K138").Value = 3233.3334
$ws.Range("L138").Value = 13034.859
$ws.Range("M138").Value = 1906.6666
$ws.Range("N138").Value = -23314.859

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3000.6667
$ws.Range("I88").Value = 3499
$ws.Range("J88").Value = 2938.375
$ws.Range("K88").Value = 3499
$ws.Range("L88").Value = 2938.375
$ws.Range("M88").Value = -3093
$ws.Range("N88").Value = -3750.375

$ws.Range("H91").Value = 3000.6667
$ws.Range("I91").Value = 3499
$ws.Range("J91").Value = 2938.375
$ws.Range("K91").Value = 3499
$ws.Range("L91").Value = 2938.375
$ws.Range("M91").Value = -2095
$ws.Range("N91").Value = -5746.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 37039828
$ws.Range("I86").Value = 58825976
$ws.Range("J86").Value = 3378.9
$ws.Range("K86").Value = 58825976
$ws.Range("L86").Value = 3378.9
$ws.Range("M86").Value = -58824853
$ws.Range("N86").Value = -5624.9

$ws.Range("H89").Value = 37039828
$ws.Range("I89").Value = 58825976
$ws.Range("J89").Value = 3378.9
$ws.Range("K89").Value = 294129880
$ws.Range("L89").Value = 16894.5
$ws.Range("M89").Value = -294124264
$ws.Range("N89").Value = -28126.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32201.912
$ws.Range("I31").Value = 2824.6296
$ws.Range("J31").Value = 145514.28
$ws.Range("K31").Value = 2824.6296
$ws.Range("L31").Value = 145514.28
$ws.Range("M31").Value = -2529.6296
$ws.Range("N31").Value = -146104.28

$ws.Range("H34").Value = 32201.912
$ws.Range("I34").Value = 2824.6296
$ws.Range("J34").Value = 145514.28
$ws.Range("K34").Value = 2824.6296
$ws.Range("L34").Value = 145514.28
$ws.Range("M34").Value = -2622.6296
$ws.Range("N34").Value = -145918.28

$ws.Range("H58").Value = 2388.0442
$ws.Range("I58").Value = 842.7907
$ws.Range("J58").Value = 5045.88
$ws.Range("K58").Value = 842.7907
$ws.Range("L58").Value = 5045.88
$ws.Range("M58").Value = -639.7907
$ws.Range("N58").Value = -5451.88

$ws.Range("H136").Value = 2388.0442
$ws.Range("I136").Value = 842.7907
$ws.Range("J136").Value = 5045.88
$ws.Range("K136").Value = 2528.3721
$ws.Range("L136").Value = 15137.64
$ws.Range("M136").Value = 21.62789999999995
$ws.Range("N136").Value = -20237.64

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 72.44444
$ws.Range("I8").Value = 72.44444
$ws.Range("K8").Value = 217.33332
$ws.Range("M8").Value = -78.33332000000001

$ws.Range("H41").Value = 270
$ws.Range("I41").Value = 10
$ws.Range("J41").Value = 400
$ws.Range("K41").Value = 30
$ws.Range("L41").Value = 1200
$ws.Range("M41").Value = 308
$ws.Range("N41").Value = -1876

$ws.Range("H69").Value = 525
$ws.Range("I69").Value = 500
$ws.Range("J69").Value = 533.3333
$ws.Range("K69").Value = 1500
$ws.Range("L69").Value = 1599.9999
$ws.Range("M69").Value = -689
$ws.Range("N69").Value = -3221.9999

$ws.Range("H70").Value = 4505.933
$ws.Range("I70").Value = 3438.9
$ws.Range("J70").Value = 6640
$ws.Range("K70").Value = 10316.7
$ws.Range("L70").Value = 19920
$ws.Range("M70").Value = -10001.7
$ws.Range("N70").Value = -20550

$ws.Range("H72").Value = 525
$ws.Range("I72").Value = 500
$ws.Range("J72").Value = 533.3333
$ws.Range("K72").Value = 4500
$ws.Range("L72").Value = 4799.9997
$ws.Range("M72").Value = -444
$ws.Range("N72").Value = -12911.9997

$ws.Range("H73").Value = 4505.933
$ws.Range("I73").Value = 3438.9
$ws.Range("J73").Value = 6640
$ws.Range("K73").Value = 10316.7
$ws.Range("L73").Value = 19920
$ws.Range("M73").Value = -9224.700000000001
$ws.Range("N73").Value = -22104

$ws.Range("H97").Value = 360.76923
$ws.Range("I97").Value = 353.33334
$ws.Range("J97").Value = 450
$ws.Range("K97").Value = 1060.00002
$ws.Range("L97").Value = 1350
$ws.Range("M97").Value = -564.0000199999999
$ws.Range("N97").Value = -2342

$ws.Range("H131").Value = 20918528
$ws.Range("I131").Value = 100202160
$ws.Range("J131").Value = 54414.527
$ws.Range("K131").Value = 300606480
$ws.Range("L131").Value = 163243.581
$ws.Range("M131").Value = -300601440
$ws.Range("N131").Value = -173323.581

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4246.9546
$ws.Range("I70").Value = 3830.6
$ws.Range("J70").Value = 4593.9165
$ws.Range("K70").Value = 3830.6
$ws.Range("L70").Value = 4593.9165
$ws.Range("M70").Value = -3560.6
$ws.Range("N70").Value = -5133.9165

$ws.Range("H73").Value = 4246.9546
$ws.Range("I73").Value = 3830.6
$ws.Range("J73").Value = 4593.9165
$ws.Range("K73").Value = 3830.6
$ws.Range("L73").Value = 4593.9165
$ws.Range("M73").Value = -2894.6
$ws.Range("N73").Value = -6465.9165

$ws.Range("H80").Value = 3239.6875
$ws.Range("I80").Value = 3233.4614
$ws.Range("J80").Value = 3266.6667
$ws.Range("K80").Value = 3233.4614
$ws.Range("L80").Value = 3266.6667
$ws.Range("M80").Value = -2235.4614
$ws.Range("N80").Value = -5262.6667

$ws.Range("H83").Value = 3239.6875
$ws.Range("I83").Value = 3233.4614
$ws.Range("J83").Value = 3266.6667
$ws.Range("K83").Value = 16167.307
$ws.Range("L83").Value = 16333.3335
$ws.Range("M83").Value = -11175.307
$ws.Range("N83").Value = -26317.3335

$ws.Range("H107").Value = 361.6842
$ws.Range("I107").Value = 256.77777
$ws.Range("J107").Value = 456.1
$ws.Range("K107").Value = 256.77777
$ws.Range("L107").Value = 456.1
$ws.Range("M107").Value = 1663.22223
$ws.Range("N107").Value = -4296.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 12470747
$ws.Range("I100").Value = 14029215
$ws.Range("K100").Value = 14029215
$ws.Range("M100").Value = -14028674

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 545.73914
$ws.Range("I107").Value = 609.93335
$ws.Range("K107").Value = 1829.80005
$ws.Range("M107").Value = 90.19994999999994

$ws.Range("H132").Value = 1425.7742
$ws.Range("I132").Value = 1200.3334
$ws.Range("J132").Value = 1899.2
$ws.Range("K132").Value = 3601.0002
$ws.Range("L132").Value = 5697.6
$ws.Range("M132").Value = -1071.0002
$ws.Range("N132").Value = -10757.6
